# Auto-generated edit script: update cryptocurrency price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '43.770.31'
$ws.Range('E2').Value = '  +4.70%  '

# Row 3
$ws.Range('D3').Value = '2.291.84'
$ws.Range('E3').Value = '  +2.99%  '

# Row 4
$ws.Range('E4').Value = '  +0.29%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.10'
$ws.Range('E5').Value = '  +0.39%  '

# Row 6
$ws.Range('E6').Value = '  +0.78%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.51'
$ws.Range('E7').Value = '  +1.44%  '

# Row 8
$ws.Range('E8').Value = '  +0.11%  '

# Row 9
$ws.Range('E9').Value = '  +4.58%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0923'
$ws.Range('E10').Value = '  +4.05%  '

# Row 11
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  +0.69%  '

# Row 12
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '2.632.60'
$ws.Range('E12').Value = '  +3.10%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.77'
$ws.Range('E13').Value = '  +0.69%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.68'
$ws.Range('E14').Value = '  +9.13%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.74'
$ws.Range('E15').Value = '  +3.31%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.812'
$ws.Range('E16').Value = '  +1.51%  '

# Row 17
$ws.Range('D17').Value = '2.288.82'
$ws.Range('E17').Value = '  +2.96%  '

# Row 18
$ws.Range('D18').Value = '43.680.67'
$ws.Range('E18').Value = '  +4.72%  '

# Row 19
$ws.Range('D19').Value = '0.0₃0931'
$ws.Range('E19').Value = '  +4.61%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.24'
$ws.Range('E20').Value = '  +0.70%  '

# Row 21
$ws.Range('E21').Value = '  +3.87%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '250.28'
$ws.Range('E22').Value = '  +0.17%  '

# Row 23
$ws.Range('E23').Value = '  -0.06%  '

# Row 24
$ws.Range('E24').Value = '  +6.98%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.36'
$ws.Range('E25').Value = '  +1.99%  '

# Row 26
$ws.Range('E26').Value = '  +2.53%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '169.87'
$ws.Range('E27').Value = '  +1.42%  '

# Row 28
$ws.Range('E28').Value = '  +0.51%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.61'
$ws.Range('E29').Value = '  +3.52%  '

# Row 30
$ws.Range('E30').Value = '  +5.58%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.65'
$ws.Range('E31').Value = '  +0.47%  '

# Row 32
$ws.Range('E32').Value = '  +0.32%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.76'
$ws.Range('E33').Value = '  +3.13%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.02'
$ws.Range('E34').Value = '  +1.51%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0659'
$ws.Range('E35').Value = '  +5.76%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.44'
$ws.Range('E36').Value = '  +3.63%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.53'
$ws.Range('E37').Value = '  -1.65%  '

# Row 38
$ws.Range('E38').Value = '  -0.70%  '

# Row 39
$ws.Range('E39').Value = '  +4.79%  '

# Row 40
$ws.Range('E40').Value = '  -0.05%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.82'
$ws.Range('E41').Value = '  +2.19%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.52'
$ws.Range('E42').Value = '  -5.30%  '

# Row 43
$ws.Range('B43').Value = 'TerraClassic'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.000218'
$ws.Range('E43').Value = '  -15.22%  '

# Row 44
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0971'
$ws.Range('E44').Value = '  -0.84%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.21'
$ws.Range('E45').Value = '  -0.07%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.60'
$ws.Range('E46').Value = '  +0.05%  '

# Row 47
$ws.Range('D47').Value = '1.473.78'
$ws.Range('E47').Value = '  +0.39%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.65'
$ws.Range('E48').Value = '  +0.87%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.29'
$ws.Range('E49').Value = '  +9.73%  '

# Row 50
$ws.Range('E50').Value = '  +1.72%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.78'
$ws.Range('E51').Value = '  -1.07%  '

# Remove temporary text-number-format overrides, keeping values as text
$forceTextCells = @('D5','D7','D10','D11','D13','D14','D15','D16','D20','D22','D25','D27','D29','D31','D33','D34','D35','D36','D37','D41','D42','D43','D44','D45','D46','D48','D49','D51')
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).ClearFormats()
}
